# Build site at 2022-09-26 16:07:08 UTC
# Rewrite the LOM3013 summary sheet: drop the long Portuguese
# paragraphs (Objetivos/Programa resumido/Programa/Bibliografia) and
# replace them with the four professor entries, reshuffle the labels
# accordingly and drop the trailing rows that are no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset all the content first so the shared-string table gets rebuilt
# from scratch in the exact order we write the cells below.
$ws.Range("A1:C27").ClearContents()

$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"
$ws.Range("B2").Value = "LOM3013"
$ws.Range("C2").Value = "LOM3013"
$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Ciência dos Materiais"
$ws.Range("C3").Value = " Ciência dos Materiais"
$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Materials Science"
$ws.Range("C4").Value = "Materials Science"
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"
$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2015"
$ws.Range("C8").Value = "01/01/2015"
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EM-3"
$ws.Range("C9").Value = "EM-3"
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Programa resumido:"
$ws.Range("B12").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C12").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("A13").Value = "Short syllabus:"
$ws.Range("A14").Value = "Programa:"
$ws.Range("B14").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C14").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("A15").Value = "Syllabus:"
$ws.Range("A16").Value = "Avaliação:"
$ws.Range("A17").Value = "Método:"
$ws.Range("B17").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C17").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("A18").Value = "Critério:"
$ws.Range("B18").Value = "Aplicação de duas provas escritas"
$ws.Range("C18").Value = "Aplicação de duas provas escritas"
$ws.Range("A19").Value = "Norma de recuperação:"
$ws.Range("B19").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2"
$ws.Range("C19").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2"
$ws.Range("A20").Value = "Bibliografia:"
$ws.Range("B20").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C20").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("A21").Value = "Requisitos:"
$ws.Range("B22").Value = "LOM3018 -  Introdução à Engenharia de Materiais  (Requisito fraco)`n"
$ws.Range("C22").Value = "LOM3018 -  Introdução à Engenharia de Materiais  (Requisito fraco)`n"

# The old sheet had five more rows (23-27); everything they used to
# hold is now folded into rows 18-22 above, so drop them.
$ws.Range("A23:A27").EntireRow.Delete()

# Re-apply the row heights used by the new layout.
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 120
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 30
